$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"
$zhCnError = "Handback file name: 4aqolpmc.oae is different with handoff file name: baa7d4ec-f1c6-4c37-a59d-e5acc889c547.12bf11a4c010425d23941690414cc9cb5fd01cb2.zh-cn."
$deDeError  = "Handback file name: 4aqolpmc.oae is different with handoff file name: baa7d4ec-f1c6-4c37-a59d-e5acc889c547.12bf11a4c010425d23941690414cc9cb5fd01cb2.de-de."

# --- Overview sheet: update status text for baa7d4ec... row (row 7) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E7").Value = $statusText
$wsOverview.Range("F7").Value = $statusText

# --- zh-cn sheet: update Status (column C) for the same row, and set Error Detail (column P) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C7").Value = $statusText
$wsZhCn.Range("P7").Value = $zhCnError
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: update Status (column C) for the same row, and set Error Detail (column P) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C7").Value = $statusText
$wsDeDe.Range("P7").Value = $deDeError
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
